$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 0.04009133333333333
$ws.Range("N2").Value = 0.120274
$ws.Range("O2").Value = 0.01033409631432067
$ws.Range("P2").Value = 0.01033409631432067
$ws.Range("Q2").Value = 0.000943977170888889
$ws.Range("R2").Value = 0.008495794538000001
$ws.Range("S2").Value = 0.00002909622060425442
$ws.Range("T2").Value = 0.00002909622060425442
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.578569084147867
$ws.Range("P3").Value = 0.578569084147867
$ws.Range("Q3").Value = 0.05284990487855556
$ws.Range("R3").Value = 0.475649143907
$ws.Range("S3").Value = 0.001628993304798166
$ws.Range("T3").Value = 0.001628993304798166
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.4110968195378122
$ws.Range("P4").Value = 0.4110968195378122
$ws.Range("Q4").Value = 0.03755200269722223
$ws.Range("R4").Value = 0.3379680242750001
$ws.Range("S4").Value = 0.001157465867083498
$ws.Range("T4").Value = 0.001157465867083498
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 0.04009133333333333
$ws.Range("N5").Value = 0.120274
$ws.Range("O5").Value = 0.01033409631432067
$ws.Range("P5").Value = 0.01033409631432067
$ws.Range("Q5").Value = 0.3308618134188889
$ws.Range("R5").Value = 2.97775632077
$ws.Range("S5").Value = 0.01019815797419616
$ws.Range("T5").Value = 0.01019815797419616
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.578569084147867
$ws.Range("P6").Value = 0.578569084147867
$ws.Range("S6").Value = 0.5709583827808372
$ws.Range("T6").Value = 0.5709583827808372
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.4110968195378122
$ws.Range("P7").Value = 0.4110968195378122
$ws.Range("S7").Value = 0.4056891072832833
$ws.Range("T7").Value = 0.4056891072832832
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 0.04009133333333333
$ws.Range("N8").Value = 0.120274
$ws.Range("O8").Value = 0.01033409631432067
$ws.Range("P8").Value = 0.01033409631432067
$ws.Range("Q8").Value = 0.003466310043777777
$ws.Range("R8").Value = 0.031196790394
$ws.Range("S8").Value = 0.000106842119520253
$ws.Range("T8").Value = 0.0001068421195202531
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.578569084147867
$ws.Range("P9").Value = 0.578569084147867
$ws.Range("S9").Value = 0.005981708062231593
$ws.Range("T9").Value = 0.005981708062231593
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.4110968195378122
$ws.Range("P10").Value = 0.4110968195378122
$ws.Range("S10").Value = 0.004250246387445456
$ws.Range("T10").Value = 0.004250246387445456
